# Update the division-problem worksheet table with a new set of problems.
# Each data row of the table (rows 1, 5, 9, 13, 17) holds five "NN÷N=" cells;
# the intervening rows are blank spacer rows. We overwrite each populated
# cell's text in place, which preserves the existing run/paragraph
# formatting (font, size, justification).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("84÷3=", "51÷5=", "30÷6=", "35÷3=", "43÷7=")
    5  = @("48÷6=", "71÷7=", "96÷4=", "65÷8=", "90÷6=")
    9  = @("38÷5=", "81÷6=", "72÷5=", "29÷3=", "38÷2=")
    13 = @("76÷2=", "74÷9=", "66÷6=", "94÷7=", "14÷5=")
    17 = @("44÷4=", "81÷3=", "46÷9=", "86÷9=", "12÷6=")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
